$d = $word.ActiveDocument

$table = $d.Tables(1)

# Row 2, Column 4 holds the "Sprint No." value ("1" -> "2")
$sprintCell = $table.Cell(2, 4)
$sprintRange = $sprintCell.Range
$sprintRange.End = $sprintRange.End - 1
$sprintRange.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "2", 2)

# Row 3, Column 2 (merged across columns 2-4) holds the Review Date ("02/09/18" -> "02/21/18")
$dateCell = $table.Cell(3, 2)
$dateRange = $dateCell.Range
$dateRange.End = $dateRange.End - 1
$dateRange.Find.Execute("02/09/18", $true, $false, $false, $false, $false, $true, 1, $false, "02/21/18", 2)
